$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update OVERALL row (row 2) and leadlag row (row 3)
# NumberFormat is forced to Text ("@") before assigning percentage-looking
# strings so Excel keeps them as literal text instead of auto-converting
# them into numeric percentage values.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C2").Value = 49
$summary.Range("D2").NumberFormat = "@"
$summary.Range("D2").Value = "65.3%"
$summary.Range("E2").NumberFormat = "@"
$summary.Range("E2").Value = "+12.2991%"
$summary.Range("F2").NumberFormat = "@"
$summary.Range("F2").Value = "+0.2510%"

$summary.Range("D3").NumberFormat = "@"
$summary.Range("D3").Value = "43.4%"
$summary.Range("E3").NumberFormat = "@"
$summary.Range("E3").Value = "+8.0975%"
$summary.Range("F3").NumberFormat = "@"
$summary.Range("F3").Value = "+0.1528%"

# ---------------------------------------------------------------------------
# Sheet "leadlag": trade #49 (row 39) transitions from OPEN to CLOSED
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")
$leadlag.Range("G39").Value = 69442.334134
$leadlag.Range("H39").Value = "CLOSED"
$leadlag.Range("I39").Value = 0.9762
$leadlag.Range("J39").Value = 9.76
$leadlag.Range("M39").Value = "time_exit_5min"
$leadlag.Range("N39").Value = 5

# ---------------------------------------------------------------------------
# Sheet "momentum": append new trade #69 (row 17), freshly OPENed
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A17").Value = 69
$momentum.Range("B17").NumberFormat = "@"
$momentum.Range("B17").Value = "2026-02-16"
$momentum.Range("C17").Value = "21:35:25"
$momentum.Range("D17").Value = "momentum"
$momentum.Range("E17").Value = "DOWN"
$momentum.Range("F17").Value = 68517.13
$momentum.Range("G17").Value = ""
$momentum.Range("H17").Value = "OPEN"
$momentum.Range("I17").Value = 0
$momentum.Range("J17").Value = 0
$momentum.Range("K17").Value = 0.9
$momentum.Range("L17").Value = "Downward momentum: -0.298% over 10 samples"
$momentum.Range("M17").Value = ""
$momentum.Range("N17").Value = 0

# ---------------------------------------------------------------------------
# Sheet "All Trades": append mirrored row 50 for the closed leadlag trade #49
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A50").Value = 49
$allTrades.Range("B50").NumberFormat = "@"
$allTrades.Range("B50").Value = "2026-02-16"
$allTrades.Range("C50").Value = "21:30:24"
$allTrades.Range("D50").Value = "leadlag"
$allTrades.Range("E50").Value = "UP"
$allTrades.Range("F50").Value = 68771.005
$allTrades.Range("G50").Value = 69442.334134
$allTrades.Range("H50").Value = "CLOSED"
$allTrades.Range("I50").Value = 0.9762
$allTrades.Range("J50").Value = 9.76
$allTrades.Range("K50").Value = 0.75
$allTrades.Range("L50").Value = "Binance leading with 0.220% move"
$allTrades.Range("M50").Value = "time_exit_5min"
$allTrades.Range("N50").Value = 5

# ---------------------------------------------------------------------------
# Sheet "Comparison": update leadlag row (row 2) statistics
# (NumberFormat forced to Text so percentage/decimal-looking strings are not
# auto-converted to numeric values by Excel.)
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")
$comparison.Range("C2").NumberFormat = "@"
$comparison.Range("C2").Value = "43.4%"
$comparison.Range("D2").NumberFormat = "@"
$comparison.Range("D2").Value = "2.63"
$comparison.Range("E2").NumberFormat = "@"
$comparison.Range("E2").Value = "+0.5683%"
$comparison.Range("G2").NumberFormat = "@"
$comparison.Range("G2").Value = "1.71"
